# "Updated users and fixtures list"
# - Final score entered for the Jul 09 match (row 50) and the Jul 10 match (row 51)
# - New fixture row added: the Jul 14, 2024 final (Spain vs England, Berlin)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scores for the two semi-final rows already on the sheet
$ws.Range("J50").Value = 2
$ws.Range("K50").Value = 1

$ws.Range("J51").Value = 1
$ws.Range("K51").Value = 2

# New fixture row 52: the final
$ws.Range("A52").Value = "Sun"

# The date column stores its values as plain text (see existing B2:B51), so
# force a text number format before writing the string to keep Excel from
# reinterpreting "Jul 14, 2024" as a date serial.
$ws.Range("B52").NumberFormat = "@"
$ws.Range("B52").Value = "Jul 14, 2024"

$ws.Range("C52").Value = "21:00:00"
$ws.Range("D52").Value = "Spain"
$ws.Range("G52").Value = "England"
$ws.Range("H52").Value = "Berlin"

# Reflect the cursor/selection position left behind in the saved file
$ws.Range("H55").Select()

Write-Output "Applied fixture list update"
